# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value2 = "Datos actualizados a 20 de Mayo de 2020 a las 13:05"

# 2) Update country case figures (row 4: Estados Unidos)
$ws.Cells.Item(4,2).Value2 = 1571131
$ws.Cells.Item(4,3).Value2 = 548
$ws.Cells.Item(4,5).Value2 = 1116346
$ws.Cells.Item(4,7).Value2 = 25
$ws.Cells.Item(4,8).Value2 = 93558

# row 13: Iran
$ws.Cells.Item(13,2).Value2 = 126949
$ws.Cells.Item(13,3).Value2 = 2346
$ws.Cells.Item(13,4).Value2 = 98808
$ws.Cells.Item(13,5).Value2 = 20958
$ws.Cells.Item(13,7).Value2 = 64
$ws.Cells.Item(13,8).Value2 = 7183

# row 28: Suiza
$ws.Cells.Item(28,2).Value2 = 30658
$ws.Cells.Item(28,3).Value2 = 40
$ws.Cells.Item(28,5).Value2 = 1067

# row 65: Oman
$ws.Cells.Item(65,5).Value2 = 4354
$ws.Cells.Item(65,7).Value2 = 1
$ws.Cells.Item(65,8).Value2 = 28

# row 126: Malta
$ws.Cells.Item(126,2).Value2 = 584
$ws.Cells.Item(126,3).Value2 = 15
$ws.Cells.Item(126,4).Value2 = 465
$ws.Cells.Item(126,5).Value2 = 113

# 3) Nepal moves ahead of Congo in the ranking (Nepal overtook Congo in
#    total cases). Swap the two countries' rows: row 133 becomes Nepal
#    (with its updated figures) and row 134 becomes Congo (figures
#    unchanged from the old row 133).
$ws.Cells.Item(133,1).Value2 = "Nepal"
$ws.Cells.Item(133,2).Value2 = 427
$ws.Cells.Item(133,3).Value2 = 25
$ws.Cells.Item(133,4).Value2 = 45
$ws.Cells.Item(133,5).Value2 = 380
$ws.Cells.Item(133,6).Value2 = 0
$ws.Cells.Item(133,7).Value2 = 0
$ws.Cells.Item(133,8).Value2 = 2

$ws.Cells.Item(134,1).Value2 = "Congo"
$ws.Cells.Item(134,2).Value2 = 420
$ws.Cells.Item(134,3).Value2 = 0
$ws.Cells.Item(134,4).Value2 = 132
$ws.Cells.Item(134,5).Value2 = 273
$ws.Cells.Item(134,6).Value2 = 0
$ws.Cells.Item(134,7).Value2 = 0
$ws.Cells.Item(134,8).Value2 = 15
